# Build site at 2023-01-09 16:18:13 UTC
#
# Fill in the "Docentes responsaveis" (responsible professors) information
# for the LOM3266 course sheet. Two new professor-name strings are
# introduced into the workbook and referenced from the cells that,
# in the published template, carry the professor information:
#   B10/C10 -> first professor (under "Objetivos:")
#   B15/C15 -> first professor again (under "Programa:")
#   B18/C18 -> second professor (under "Metodo:")
# Row 13 (under "Programa resumido:") keeps carrying the activation
# date string, same text as B8/C8.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$professor1 = "519033 - Carlos Yujiro Shigue"
$professor2 = "7290967 - Emerson Gonçalves de Melo"
$ativacao   = "01/01/2023"

# --- B10 / C10 -------------------------------------------------------
$ws.Range("B10").Value = $professor1
$ws.Range("C10").Value = $professor1

# --- B13 / C13 ---------------------------------------------------------
# "01/01/2023" looks like a date, so a plain Value assignment on a
# General-formatted cell would be auto-converted into a date serial
# number by Excel. Write it on a throw-away cell that has been forced
# to Text format first, then copy *values only* (not formats) onto the
# real target cells so their own number format / style stays untouched.
$scratch = $ws.Range("ZZ1000")
$scratch.NumberFormat = "@"
$scratch.Value = $ativacao
$scratch.Copy()
$ws.Range("B13").PasteSpecial(-4163)  # xlPasteValues
$ws.Range("C13").PasteSpecial(-4163)  # xlPasteValues
$scratch.Clear()
$excel.CutCopyMode = $false

# --- B15 / C15 -------------------------------------------------------
$ws.Range("B15").Value = $professor1
$ws.Range("C15").Value = $professor1

# --- B18 / C18 -------------------------------------------------------
$ws.Range("B18").Value = $professor2
$ws.Range("C18").Value = $professor2
